$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 551.5714
$ws.Range("I2").Value = 563.3333
$ws.Range("J2").Value = 542.75
$ws.Range("K2").Value = 563.3333
$ws.Range("L2").Value = 542.75
$ws.Range("M2").Value = -450.3333
$ws.Range("N2").Value = -768.75
$ws.Range("H40").Value = 7324.75
$ws.Range("J40").Value = 2449
$ws.Range("L40").Value = 2449
$ws.Range("N40").Value = -2799
$ws.Range("H112").Value = 2609.2974
$ws.Range("I112").Value = 1494.2858
$ws.Range("J112").Value = 2869.4666
$ws.Range("K112").Value = 4482.857400000001
$ws.Range("L112").Value = 8608.399800000001
$ws.Range("M112").Value = -3374.857400000001
$ws.Range("N112").Value = -10824.3998
$ws.Range("H125").Value = 949.8
$ws.Range("I125").Value = 956.125
$ws.Range("J125").Value = 924.5
$ws.Range("K125").Value = 8605.125
$ws.Range("L125").Value = 8320.5
$ws.Range("M125").Value = -6145.125
$ws.Range("N125").Value = -13240.5
$ws.Range("H138").Value = 2928.0188
$ws.Range("J138").Value = 3005.6511
$ws.Range("L138").Value = 9016.953300000001
$ws.Range("N138").Value = -19296.9533

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 324.27274
$ws.Range("I4").Value = 317.8
$ws.Range("J4").Value = 389
$ws.Range("K4").Value = 317.8
$ws.Range("L4").Value = 389
$ws.Range("M4").Value = -201.8
$ws.Range("N4").Value = -621
$ws.Range("H32").Value = 18013.525
$ws.Range("I32").Value = 6500.727
$ws.Range("K32").Value = 6500.727
$ws.Range("M32").Value = -6213.727
$ws.Range("H54").Value = 8000
$ws.Range("I54").Value = 8000
$ws.Range("K54").Value = 8000
$ws.Range("M54").Value = -7231
$ws.Range("H97").Value = 1219.6923
$ws.Range("I97").Value = 1282.8182
$ws.Range("K97").Value = 1282.8182
$ws.Range("M97").Value = -786.8181999999999
$ws.Range("H102").Value = 2598.45
$ws.Range("I102").Value = 2628.2942
$ws.Range("J102").Value = 2429.3333
$ws.Range("K102").Value = 2628.2942
$ws.Range("L102").Value = 2429.3333
$ws.Range("M102").Value = -1006.2942
$ws.Range("N102").Value = -5673.3333

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 569.4286
$ws.Range("J22").Value = 500
$ws.Range("L22").Value = 500
$ws.Range("N22").Value = -846
$ws.Range("H76").Value = 26748.25
$ws.Range("J76").Value = 26748.25
$ws.Range("L76").Value = 26748.25
$ws.Range("N76").Value = -27378.25
$ws.Range("H79").Value = 26748.25
$ws.Range("J79").Value = 26748.25
$ws.Range("L79").Value = 26748.25
$ws.Range("N79").Value = -28932.25
$ws.Range("H107").Value = 874.12
$ws.Range("I107").Value = 950.8823
$ws.Range("K107").Value = 950.8823
$ws.Range("M107").Value = 969.1177

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 125134.625
$ws.Range("I7").Value = 166804.5
$ws.Range("K7").Value = 166804.5
$ws.Range("M7").Value = -166691.5
$ws.Range("H15").Value = 5240.5
$ws.Range("I15").Value = 1008
$ws.Range("J15").Value = 9473
$ws.Range("K15").Value = 1008
$ws.Range("L15").Value = 9473
$ws.Range("M15").Value = -838
$ws.Range("N15").Value = -9813

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 473.3
$ws.Range("I5").Value = 422
$ws.Range("K5").Value = 1266
$ws.Range("M5").Value = -1154
$ws.Range("H11").Value = 942.4286
$ws.Range("I11").Value = 942.4286
$ws.Range("K11").Value = 2827.2858
$ws.Range("M11").Value = -2687.2858
$ws.Range("H81").Value = 55406.5
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 55406.5
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H88").Value = 4527.25
$ws.Range("J88").Value = 5136.6665
$ws.Range("L88").Value = 15409.9995
$ws.Range("N88").Value = -16265.9995
$ws.Range("H91").Value = 4527.25
$ws.Range("J91").Value = 5136.6665
$ws.Range("L91").Value = 15409.9995
$ws.Range("N91").Value = -18373.9995
$ws.Range("H97").Value = 1999.5
$ws.Range("J97").Value = 1999.5
$ws.Range("L97").Value = 5998.5
$ws.Range("N97").Value = -6990.5
$ws.Range("H102").Value = 2600
$ws.Range("I102").Value = 2600
$ws.Range("K102").Value = 7800
$ws.Range("M102").Value = -5366
$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("M103").ClearContents()
$ws.Range("N103").ClearContents()
$ws.Range("H107").Value = 596.05554
$ws.Range("I107").Value = 309
$ws.Range("J107").Value = 825.7
$ws.Range("K107").Value = 927
$ws.Range("L107").Value = 2477.1
$ws.Range("M107").Value = 993
$ws.Range("N107").Value = -6317.1
$ws.Range("H108").Value = 1024
$ws.Range("I108").Value = 1024
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 3072
$ws.Range("L108").Value = 0
$ws.Range("M108").Value = -192
$ws.Range("N108").ClearContents()
$ws.Range("H113").Value = 748.75
$ws.Range("I113").Value = 295.25
$ws.Range("K113").Value = 885.75
$ws.Range("M113").Value = 1284.25
$ws.Range("H122").Value = 1730.2084
$ws.Range("J122").Value = 1864.6875
$ws.Range("L122").Value = 16782.1875
$ws.Range("N122").Value = -21682.1875
$ws.Range("H125").Value = 9799.833000000001
$ws.Range("I125").Value = 9799
$ws.Range("K125").Value = 29397
$ws.Range("M125").Value = -24477
$ws.Range("H135").Value = 473.3
$ws.Range("I135").Value = 422
$ws.Range("K135").Value = 3798
$ws.Range("M135").Value = -1263

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 793.25
$ws.Range("J97").Value = 521.2
$ws.Range("L97").Value = 521.2
$ws.Range("N97").Value = -1513.2
$ws.Range("H102").Value = 3185.4375
$ws.Range("I102").Value = 3045.8
$ws.Range("K102").Value = 3045.8
$ws.Range("M102").Value = -1423.8
$ws.Range("H113").Value = 2257.4707
$ws.Range("I113").Value = 2125.2
$ws.Range("K113").Value = 2125.2
$ws.Range("M113").Value = 44.80000000000018
$ws.Range("H122").Value = 1693.5714
$ws.Range("I122").Value = 1148.6
$ws.Range("K122").Value = 3445.8
$ws.Range("M122").Value = -995.7999999999997
$ws.Range("H126").Value = 3901.8333
$ws.Range("I126").Value = 4139
$ws.Range("J126").Value = 3732.4285
$ws.Range("K126").Value = 12417
$ws.Range("L126").Value = 11197.2855
$ws.Range("M126").Value = -9947
$ws.Range("N126").Value = -16137.2855

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7185.8423
$ws.Range("I7").Value = 6477.25
$ws.Range("J7").Value = 10965
$ws.Range("K7").Value = 6477.25
$ws.Range("L7").Value = 10965
$ws.Range("M7").Value = -6365.25
$ws.Range("N7").Value = -11189
$ws.Range("H55").Value = 409.3
$ws.Range("I55").Value = 596.1667
$ws.Range("K55").Value = 596.1667
$ws.Range("M55").Value = -423.1667
$ws.Range("H62").Value = 22997.25
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H64").Value = 39535.5
$ws.Range("J64").Value = 39535.5
$ws.Range("L64").Value = 39535.5
$ws.Range("N64").Value = -39985.5
$ws.Range("H65").Value = 22997.25
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H67").Value = 39535.5
$ws.Range("J67").Value = 39535.5
$ws.Range("L67").Value = 39535.5
$ws.Range("N67").Value = -41095.5
$ws.Range("H82").Value = 2087.5454
$ws.Range("I82").Value = 2087.5454
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 2087.5454
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -1726.5454
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 2087.5454
$ws.Range("I85").Value = 2087.5454
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 2087.5454
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -839.5454
$ws.Range("N85").ClearContents()
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H126").Value = 7185.8423
$ws.Range("I126").Value = 6477.25
$ws.Range("J126").Value = 10965
$ws.Range("K126").Value = 19431.75
$ws.Range("L126").Value = 32895
$ws.Range("M126").Value = -16961.75
$ws.Range("N126").Value = -37835

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 44208.168
$ws.Range("J63").Value = 44208.168
$ws.Range("L63").Value = 44208.168
$ws.Range("N63").Value = -45456.168
$ws.Range("H66").Value = 44208.168
$ws.Range("J66").Value = 44208.168
$ws.Range("L66").Value = 132624.504
$ws.Range("N66").Value = -138864.504
$ws.Range("H126").Value = 6585.4287
$ws.Range("I126").Value = 4719.8
$ws.Range("K126").Value = 14159.4
$ws.Range("M126").Value = -11689.4
$ws.Range("H127").Value = 53750
$ws.Range("J127").Value = 53750
$ws.Range("L127").Value = 53750
$ws.Range("N127").Value = -63670
$ws.Range("H132").Value = 2332.5
$ws.Range("I132").Value = 2402.6843
$ws.Range("J132").Value = 999
$ws.Range("K132").Value = 7208.0529
$ws.Range("L132").Value = 2997
$ws.Range("M132").Value = -4678.0529
$ws.Range("N132").Value = -8057

Write-Output "Applied 249 cell changes across 8 sheets"